# "fall 22 week 8 complete" - append week 8 matchup results (24 new rows)
# to Sheet1, directly below the existing data (previously ending at row 1337).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 1338
$endRow = 1361
$data = New-Object 'object[,]' 24,4
$data[0,0] = 4
$data[0,1] = 0
$data[0,2] = 4
$data[0,3] = 2
$data[1,0] = 5
$data[1,1] = 2
$data[1,2] = 5
$data[1,3] = 1
$data[2,0] = 3
$data[2,1] = 2
$data[2,2] = 4
$data[2,3] = 1
$data[3,0] = 3
$data[3,1] = 2
$data[3,2] = 3
$data[3,3] = 1
$data[4,0] = 7
$data[4,1] = 2
$data[4,2] = 5
$data[4,3] = 0
$data[5,0] = 5
$data[5,1] = 2
$data[5,2] = 4
$data[5,3] = 1
$data[6,0] = 4
$data[6,1] = 3
$data[6,2] = 3
$data[6,3] = 0
$data[7,0] = 4
$data[7,1] = 1
$data[7,2] = 6
$data[7,3] = 2
$data[8,0] = 6
$data[8,1] = 2
$data[8,2] = 7
$data[8,3] = 1
$data[9,0] = 2
$data[9,1] = 1
$data[9,2] = 3
$data[9,3] = 2
$data[10,0] = 3
$data[10,1] = 1
$data[10,2] = 4
$data[10,3] = 2
$data[11,0] = 6
$data[11,1] = 0
$data[11,2] = 6
$data[11,3] = 2
$data[12,0] = 5
$data[12,1] = 2
$data[12,2] = 5
$data[12,3] = 0
$data[13,0] = 5
$data[13,1] = 0
$data[13,2] = 4
$data[13,3] = 2
$data[14,0] = 4
$data[14,1] = 2
$data[14,2] = 4
$data[14,3] = 1
$data[15,0] = 5
$data[15,1] = 3
$data[15,2] = 4
$data[15,3] = 0
$data[16,0] = 5
$data[16,1] = 2
$data[16,2] = 5
$data[16,3] = 0
$data[17,0] = 6
$data[17,1] = 1
$data[17,2] = 6
$data[17,3] = 2
$data[18,0] = 3
$data[18,1] = 2
$data[18,2] = 4
$data[18,3] = 0
$data[19,0] = 4
$data[19,1] = 1
$data[19,2] = 4
$data[19,3] = 2
$data[20,0] = 6
$data[20,1] = 0
$data[20,2] = 5
$data[20,3] = 2
$data[21,0] = 6
$data[21,1] = 3
$data[21,2] = 5
$data[21,3] = 0
$data[22,0] = 3
$data[22,1] = 2
$data[22,2] = 5
$data[22,3] = 0
$data[23,0] = 5
$data[23,1] = 0
$data[23,2] = 5
$data[23,3] = 2

$rng = $ws.Range("A" + $startRow + ":D" + $endRow)
$rng.Value = $data

# Scroll the view down near the newly-added rows and move the active
# selection to the cell just below the new data (matches the workbook's
# saved cursor position after the edit).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1349
$ws.Range("A" + ($endRow + 1)).Select()

